# Actualiza cronograma con avances.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 70
$ws.Range("E8").Value = 60
$ws.Range("E21").Value = 50
